$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experimentos")

# Fill in new values for row 4 (G4, H4)
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = 0.17740429505135399

# Update the selection to H7 (single cell)
$ws.Activate()
$ws.Range("H7").Select()
